$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Notes master "cached" date field bump: 17/09/2023 -> 18/09/2023
# ---------------------------------------------------------------------------
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "18/09/2023"

# ---------------------------------------------------------------------------
# 2) Slide 1 - "Go to www.menti.com and use the code 4831 7794" textbox
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(7)

# Move/resize the textbox.
$shp.Left = 634.6279296875
$shp.Top = 440.46392822265625
$shp.Width = 282.25482177734375
$shp.Height = 50.892208099365234

$tr = $shp.TextFrame.TextRange

# Work from the end of the text backwards so earlier character offsets stay
# valid while later pieces are edited.

# "4831 7794" -> " 1727 7990"
$tr.Characters(39, 9).Text = " 1727 7990"

# "use the code " -> "with the code " (adds the new "with" word, drops "use")
$tr.Characters(26, 13).Text = "with the code "
# Split "with" into its own run (clean formatting, no bold/italic baggage)
$tr.Characters(26, 4).Font.Name = "Montserrat"

# " and " -> " "
$tr.Characters(20, 5).Text = " "

# "www.menti.com" -> "menti.com"
$tr.Characters(7, 13).Text = "menti.com"

# "Go to " -> "Ask questions on " (adds the new "questions on " phrase)
$tr.Characters(1, 6).Text = "Ask questions on "
# Split "questions on " into its own run (clean formatting, no bold/italic baggage)
$tr.Characters(5, 13).Font.Name = "Montserrat"
